$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: delete row 264 (004452912 BRUNO -19464), old trailing BRUNO entry ---
$ws.Rows(264).Delete()

# --- Step 2: delete row 206 (004472431 LUIS 33.08), old LUIS entry ---
$ws.Rows(206).Delete()

# --- Step 3: delete row 198 (003894173 ANDREA, old position after SILVANIA) ---
$ws.Rows(198).Delete()

# --- Step 4: delete row 197 (005070742 JUSCELINO, old position after SILVANIA) ---
$ws.Rows(197).Delete()

# --- Step 5: update row 196 (004207955 SILVANIA) value 44.54 -> 42.32 ---
$ws.Cells.Item(196,3).Value = 42.32

# --- Step 6: insert JUSCELINO + ANDREA rows before row 196 (SILVANIA) ---
$ws.Rows(196).Insert()
$ws.Cells.Item(196,1).NumberFormat = "@"
$ws.Cells.Item(196,1).Value = "003894173"
$ws.Cells.Item(196,2).Value = "ANDREA"
$ws.Cells.Item(196,3).Value = 44.01
$ws.Rows(196).Insert()
$ws.Cells.Item(196,1).NumberFormat = "@"
$ws.Cells.Item(196,1).Value = "005070742"
$ws.Cells.Item(196,2).Value = "JUSCELINO"
$ws.Cells.Item(196,3).Value = 44.06

# --- Step 7: delete row 189 (001761119 BLUEMETRIX) ---
$ws.Rows(189).Delete()

# --- Step 8: insert ROGERIO, ELENI, VALMIR rows before row 10 (005152037 RODRIGO) ---
$ws.Rows(10).Insert()
$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = "004487140"
$ws.Cells.Item(10,2).Value = "VALMIR"
$ws.Cells.Item(10,3).Value = 9841.47
$ws.Rows(10).Insert()
$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = "005101676"
$ws.Cells.Item(10,2).Value = "ELENI"
$ws.Cells.Item(10,3).Value = 10290.89
$ws.Rows(10).Insert()
$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = "004487016"
$ws.Cells.Item(10,2).Value = "ROGERIO"
$ws.Cells.Item(10,3).Value = 10862.74

# --- Step 9: insert HEPTA, RENATA, LUIS(39177.22), EDINARDO rows before row 6 (004224011 THOMAS) ---
$ws.Rows(6).Insert()
$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "004211807"
$ws.Cells.Item(6,2).Value = "EDINARDO"
$ws.Cells.Item(6,3).Value = 29500
$ws.Rows(6).Insert()
$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "004472431"
$ws.Cells.Item(6,2).Value = "LUIS"
$ws.Cells.Item(6,3).Value = 39177.22
$ws.Rows(6).Insert()
$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "004646727"
$ws.Cells.Item(6,2).Value = "RENATA"
$ws.Cells.Item(6,3).Value = 44546.09
$ws.Rows(6).Insert()
$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "004359408"
$ws.Cells.Item(6,2).Value = "HEPTA"
$ws.Cells.Item(6,3).Value = 53372.89

# --- Step 10: insert BRUNO(106707.28) row before row 3 (005081833 PEDRO) ---
$ws.Rows(3).Insert()
$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "004452912"
$ws.Cells.Item(3,2).Value = "BRUNO"
$ws.Cells.Item(3,3).Value = 106707.28
